$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Multiply the "particip" (E) and "taxa_sucesso" (F) columns by 100 for
# the data rows (2-7), as those values were stored as fractions (0-1)
# but should be stored as whole percentage numbers (0-100).
for ($r = 2; $r -le 7; $r++) {
    $eCell = $ws.Cells.Item($r, 5)  # Column E
    $fCell = $ws.Cells.Item($r, 6)  # Column F

    $eCell.Value = $eCell.Value2 * 100
    $fCell.Value = $fCell.Value2 * 100
}
